# Add the three new submission names to sheet "14" and make it the
# active/selected sheet, matching the author's "updating names on excel
# file" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("14")

$ws.Range("A1").Value = "Alon Nadel"
$ws.Range("A2").Value = "Bar Naor"
$ws.Range("A3").Value = "Gil Tamir"

# Make this the active tab/sheet and leave the selection on the next
# empty row, just like the author left it.
$ws.Activate()
[void]$ws.Range("A4").Select()
